$d = $word.ActiveDocument

# Helper: find $anchor, then within the range that follows it (to the end of
# the story) replace the first occurrence of $pattern with "" (i.e. delete
# it). Restricting the search to start right after the anchor keeps the
# anchor's own run intact instead of merging it with what follows.
function Remove-TextAfter([string]$anchor, [string]$pattern) {
    $anchorRng = $d.Content
    $anchorRng.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $afterAnchor = $anchorRng.End
    $scoped = $d.Range($afterAnchor, $d.Content.End)
    $scoped.Find.Execute($pattern, $false, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null
}

# 1) Drop the "Course GPA: 4.00" line from the CS1301xI entry.
Remove-TextAfter "11/5/18 – 11/18/18" "^lCourse GPA: 4.00"

# 2) Drop the "Course GPA: 4.00" line from the CS1301xII entry.
Remove-TextAfter "1/26/18" "^lCourse GPA: 4.00"

# 3) Drop the "Course GPA: 4.00" line (plus the trailing break+tab that used
#    to precede the _GoBack bookmark) from the Univ. of Colorado entry.
Remove-TextAfter "04/10/19" "^lCourse GPA: 4.00^l^t"

# 4) The _GoBack bookmark used to sit at the very end of the document; move
#    it into the CS1301xI course title, splitting "...Computing in" from
#    " Python I...". Delete the old one first, then re-add at the new spot.
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$splitRng = $d.Content
$splitRng.Find.Execute("CS1301xI – Computing in", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $d.Range($splitRng.End, $splitRng.End)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# 5) Append the new "University of Adelaide" paragraph right after the
#    "Requirements Gathering..." (Univ. of Colorado) paragraph.
$coloradoCourse = $d.Paragraphs(5)
$insertionPoint = $d.Range($coloradoCourse.Range.End, $coloradoCourse.Range.End)
$newParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>University of Adelaide</w:t></w:r><w:r><w:br/></w:r><w:r><w:tab/></w:r><w:r><w:t>Think. Create. Code.</w:t></w:r><w:r><w:br/></w:r><w:r><w:tab/><w:t>05/23/19 – present</w:t></w:r></w:p></w:body></w:document>'
$insertionPoint.InsertXML($newParagraphXml)

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
